# Hortaliza / Pepino ensalada - weekly price update
# Insert two new weekly observation rows near the top of the price
# history block (row 351) and two new rows at the bottom (474/475),
# shifting the existing rows 351..473 down to 353..475.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert two blank rows at row 351, pushing old rows 351-473 to 353-475
$ws.Range("A351:A352").EntireRow.Insert()

# --- 2) Fill the two new rows (351, 352) with the new weekly observations
# Row 351 - Primera
$ws.Range("A351").Value = 1
$ws.Range("B351").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C351").Value = "Arica y Parinacota"
$ws.Range("D351").Value = 45027
$ws.Range("E351").Value = 15
$ws.Range("F351").Value = 100112043
$ws.Range("G351").Value = "Pepino ensalada"
$ws.Range("H351").Value = "Sin especificar"
$ws.Range("I351").Value = "Primera"
$ws.Range("J351").Value = 220
$ws.Range("K351").Value = 4000
$ws.Range("L351").Value = 5000
$ws.Range("M351").Value = 4455
$ws.Range("N351").Value = "$/caja 70 unidades"
$ws.Range("O351").Value = "Región de Arica y Parinacota"
$ws.Range("P351").Value = 64
$ws.Range("Q351").Value = 70
$ws.Range("R351").Value = "Hortaliza"

# Row 352 - Segunda
$ws.Range("A352").Value = 1
$ws.Range("B352").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C352").Value = "Arica y Parinacota"
$ws.Range("D352").Value = 45027
$ws.Range("E352").Value = 15
$ws.Range("F352").Value = 100112043
$ws.Range("G352").Value = "Pepino ensalada"
$ws.Range("H352").Value = "Sin especificar"
$ws.Range("I352").Value = "Segunda"
$ws.Range("J352").Value = 250
$ws.Range("K352").Value = 3000
$ws.Range("L352").Value = 4000
$ws.Range("M352").Value = 3600
$ws.Range("N352").Value = "$/caja 100 unidades"
$ws.Range("O352").Value = "Región de Arica y Parinacota"
$ws.Range("P352").Value = 36
$ws.Range("Q352").Value = 100
$ws.Range("R352").Value = "Hortaliza"

# --- 3) Append two new rows (474, 475) at the bottom of the table
# Row 474 - Primera
$ws.Range("A474").Value = 1
$ws.Range("B474").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C474").Value = "Arica y Parinacota"
$ws.Range("D474").Value = 45007
$ws.Range("E474").Value = 15
$ws.Range("F474").Value = 100112043
$ws.Range("G474").Value = "Pepino ensalada"
$ws.Range("H474").Value = "Sin especificar"
$ws.Range("I474").Value = "Primera"
$ws.Range("J474").Value = 130
$ws.Range("K474").Value = 3500
$ws.Range("L474").Value = 4000
$ws.Range("M474").Value = 3750
$ws.Range("N474").Value = "$/caja 70 unidades"
$ws.Range("O474").Value = "Región de Arica y Parinacota"
$ws.Range("P474").Value = 54
$ws.Range("Q474").Value = 70
$ws.Range("R474").Value = "Hortaliza"

# Row 475 - Segunda
$ws.Range("A475").Value = 1
$ws.Range("B475").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C475").Value = "Arica y Parinacota"
$ws.Range("D475").Value = 45007
$ws.Range("E475").Value = 15
$ws.Range("F475").Value = 100112043
$ws.Range("G475").Value = "Pepino ensalada"
$ws.Range("H475").Value = "Sin especificar"
$ws.Range("I475").Value = "Segunda"
$ws.Range("J475").Value = 150
$ws.Range("K475").Value = 3000
$ws.Range("L475").Value = 3500
$ws.Range("M475").Value = 3250
$ws.Range("N475").Value = "$/caja 100 unidades"
$ws.Range("O475").Value = "Región de Arica y Parinacota"
$ws.Range("P475").Value = 32
$ws.Range("Q475").Value = 100
$ws.Range("R475").Value = "Hortaliza"
